$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28..53 down to 29..54
$ws.Rows(28).Insert()

# Populate the newly inserted row 28 with the new weekly record
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44596
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112027
$ws.Range("G28").Value = "Melón"
$ws.Range("H28").Value = "Tuna"
$ws.Range("I28").Value = "Segunda"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = 5500
$ws.Range("N28").Value = "$/caja 24 unidades"
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 229
$ws.Range("Q28").Value = 24
$ws.Range("R28").Value = "Hortaliza"
